$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reorders (permutes) the data rows 2-10: the content that used to
# live in one row now lives in a different row (the header row 1 and the
# columns are unchanged - this is purely a row shuffle). We move data with
# Range.Copy (rather than writing scalar .Value/.Value2) so each cell keeps
# its original stored type - otherwise numeric-looking text such as
# "2023-08-12" or the digit-only strings in the "Antal" column would get
# silently reinterpreted as a real date/number.

$lastColIndex = 51   # column AY
$srcFirst = 2
$srcLast = 10
$stagingOffset = 200

# new worksheet row -> row whose data should end up there
$map = @{
    2  = 9
    3  = 4
    4  = 5
    5  = 8
    6  = 10
    7  = 2
    8  = 6
    9  = 3
    10 = 7
}

# Step 1: copy each existing data row into a staging area far below the
# used range so source and destination never overlap while we rearrange
# them row-by-row.
for ($r = $srcFirst; $r -le $srcLast; $r++) {
    $srcRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $lastColIndex))
    $stageRow = $r + $stagingOffset
    $dstRange = $ws.Range($ws.Cells.Item($stageRow, 1), $ws.Cells.Item($stageRow, $lastColIndex))
    $srcRange.Copy($dstRange)
}

# Step 2: blank out the original rows. Copy() in this runtime only
# overwrites a destination cell when the source cell actually holds
# something - a blank source cell left the previous destination content
# untouched - so without this the rows would keep stray leftovers from
# whatever used to sit in that row position.
$ws.Range($ws.Cells.Item($srcFirst, 1), $ws.Cells.Item($srcLast, $lastColIndex)).ClearContents()

# Step 3: copy from the staged copies back into the final row positions in
# the permuted order, one cell at a time so only columns that actually hold
# data (including legitimate empty-string cells) are (re)written - this
# keeps untouched/always-blank columns genuinely absent, matching the
# original sparse layout instead of materialising empty cells everywhere.
foreach ($newRow in ($map.Keys | Sort-Object)) {
    $oldRow = $map[$newRow]
    $stageRow = $oldRow + $stagingOffset
    for ($c = 1; $c -le $lastColIndex; $c++) {
        $stageCell = $ws.Cells.Item($stageRow, $c)
        if ($stageCell.Value2 -ne $null) {
            $destCell = $ws.Cells.Item($newRow, $c)
            $stageCell.Copy($destCell)
        }
    }
}

# Step 4: clear the staging area so the used range collapses back down to
# its original A1:AY10 extent.
$stageClearFirst = $srcFirst + $stagingOffset
$stageClearLast = $srcLast + $stagingOffset
$ws.Range($ws.Cells.Item($stageClearFirst, 1), $ws.Cells.Item($stageClearLast, $lastColIndex)).Clear()
